# Updated cryptos list on Fri Sep 15 06:11:19 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Some "Price" values in column D are plain decimal numbers (e.g. "213.25").
    # Excel auto-converts a bare numeric string typed into a cell into a Number,
    # but these columns must stay Text (others contain thousand-separated
    # values like "26.652.95" that can't be numeric anyway). Force Text first,
    # write the value, then drop back to the Normal style so no lingering
    # number-format style is left on the cell.
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.652.95"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.633.77"
$ws.Range("E3").Value = "  +0.95%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "213.25"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6 - XRP
Set-TextValue "D6" "0.503"
$ws.Range("E6").Value = "  +3.81%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.09%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +2.48%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0623"

# Row 10 - Solana
Set-TextValue "D10" "19.24"
$ws.Range("E10").Value = "  +2.80%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0845"
$ws.Range("E11").Value = "  +3.49%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.862.63"
$ws.Range("E12").Value = "  +1.07%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.666.13"
$ws.Range("E13").Value = "  +2.97%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +2.59%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.86%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.644.29"
$ws.Range("E16").Value = "  +1.36%  "

# Row 17 - Litecoin
Set-TextValue "D17" "63.49"
$ws.Range("E17").Value = "  +2.06%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +2.51%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "218.80"
$ws.Range("E19").Value = "  +8.73%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.02%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.59%  "

# Row 22 & 23 - Chainlink and Avalanche swap ranking positions
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D22" "9.46"
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D23" "6.22"
$ws.Range("E23").Value = "  +2.91%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +2.10%  "

# Row 25 - Monero
Set-TextValue "D25" "148.52"
$ws.Range("E25").Value = "  +2.84%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.11%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.36%  "

# Row 28 - Cosmos
Set-TextValue "D28" "6.92"
$ws.Range("E28").Value = "  +5.91%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.50"
$ws.Range("E29").Value = "  +2.36%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.58%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.23%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.31"
$ws.Range("E32").Value = "  +4.38%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +2.90%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.51"
$ws.Range("E34").Value = "  +1.45%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.04%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.208.36"
$ws.Range("E36").Value = "  +2.74%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.0172"
$ws.Range("E37").Value = "  +5.38%  "

# Row 38 - ARBITRUM
Set-TextValue "D38" "0.811"
$ws.Range("E38").Value = "  +0.87%  "

# Row 39 - PaxDollar
$ws.Range("E39").Value = "  +0.07%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  +2.14%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  -1.08%  "

# Row 42 - FraxShare
Set-TextValue "D42" "5.44"
$ws.Range("E42").Value = "  +1.87%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "0.792"
$ws.Range("E43").Value = "  +0.48%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.770.96"
$ws.Range("E44").Value = "  +0.96%  "

# Row 45 - Quant
Set-TextValue "D45" "93.26"
$ws.Range("E45").Value = "  +0.84%  "

# Row 46 - RenderToken
Set-TextValue "D46" "1.56"
$ws.Range("E46").Value = "  +1.97%  "

# Row 47 - Aave
Set-TextValue "D47" "54.81"
$ws.Range("E47").Value = "  +2.26%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +1.06%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "7.73"
$ws.Range("E49").Value = "  +6.28%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  +0.51%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  +0.30%  "
